# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.759.19'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').Value = '2.109.68'
$ws.Range('E3').Value = '  +5.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.84'
$ws.Range('E5').Value = '  +3.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5295'
$ws.Range('E7').Value = '  +3.67%  '
$ws.Range('E8').Value = '  +5.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08960'
$ws.Range('E9').Value = '  +2.54%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.21'
$ws.Range('E10').Value = '  +10.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.168'
$ws.Range('E11').Value = '  +2.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.78'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '2.104.42'
$ws.Range('E13').Value = '  +5.26%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.726'
$ws.Range('E14').Value = '  +2.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.766'
$ws.Range('E15').Value = '  +4.15%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '96.73'
$ws.Range('E16').Value = '  +2.65%  '
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('E18').Value = '  +1.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06689'
$ws.Range('E19').Value = '  +2.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.00'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.302'
$ws.Range('E22').Value = '  +2.65%  '
$ws.Range('D23').Value = '30.812.46'
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.27'
$ws.Range('E24').Value = '  +4.45%  '
$ws.Range('D25').Value = '2.352.66'
$ws.Range('E25').Value = '  +5.40%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.283'
$ws.Range('E26').Value = '  +2.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '22.60'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.570'
$ws.Range('E28').Value = '  +6.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '162.44'
$ws.Range('E29').Value = '  -0.88%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.17'
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.191'
$ws.Range('E31').Value = '  +4.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1081'
$ws.Range('E32').Value = '  +2.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.170'
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.958'
$ws.Range('E34').Value = '  +3.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.548'
$ws.Range('E35').Value = '  +15.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02597'
$ws.Range('E36').Value = '  +3.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '9.624'
$ws.Range('E37').Value = '  +7.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.547'
$ws.Range('E38').Value = '  +3.08%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06754'
$ws.Range('E39').Value = '  +2.49%  '
$ws.Range('E40').Value = '  +3.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2272'
$ws.Range('E41').Value = '  +3.14%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6820'
$ws.Range('E42').Value = '  +2.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.247'
$ws.Range('E43').Value = '  +1.77%  '
$ws.Range('E44').Value = '  +0.06%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '14.08'
$ws.Range('E45').Value = '  +3.62%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6412'
$ws.Range('E46').Value = '  +4.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.224'
$ws.Range('E47').Value = '  +1.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.654'
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.261'
$ws.Range('E49').Value = '  -0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '83.04'
$ws.Range('E50').Value = '  +3.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.194'
$ws.Range('E51').Value = '  +8.06%  '
